$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.870.65'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '3.116.57'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.98'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.92'
$ws.Range('E6').Value = '  +1.99%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.521'
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('B9').Value = 'Toncoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.44'
$ws.Range('E9').Value = '  -2.60%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.153'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.484'
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000246'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.98'
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.123'
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.632.01'
$ws.Range('E15').Value = '  +0.93%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '66.833.61'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.18'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.114.09'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.30'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '477.95'
$ws.Range('E20').Value = '  +2.65%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.98'
$ws.Range('E21').Value = '  +6.01%  '
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '84.08'
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.33'
$ws.Range('E24').Value = '  +3.42%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.28'
$ws.Range('E25').Value = '  -3.32%  '
$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.07'
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('B28').Value = 'NEARProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.87'
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.37'
$ws.Range('E29').Value = '  -0.98%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.67'
$ws.Range('E30').Value = '  +0.40%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '28.57'
$ws.Range('E31').Value = '  +1.50%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.115'
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D33').Value = '0.0₃0941'
$ws.Range('E33').Value = '  -7.39%  '
$ws.Range('B34').Value = 'FirstDigitalUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.85'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('B36').Value = 'Mantle'
$ws.Range('C36').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.972'
$ws.Range('E36').Value = '  -2.94%  '
$ws.Range('B37').Value = 'Arweave'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '47.22'
$ws.Range('E37').Value = '  +0.98%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.19'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.06'
$ws.Range('E39').Value = '  -3.18%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.310'
$ws.Range('E40').Value = '  -2.55%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.122'
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.67'
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.823.54'
$ws.Range('E43').Value = '  +2.13%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '383.08'
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0354'
$ws.Range('E45').Value = '  -1.63%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.56'
$ws.Range('E46').Value = '  -9.03%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '135.01'
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.78'
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.20'
$ws.Range('E50').Value = '  -1.81%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.108'
$ws.Range('E51').Value = '  -0.66%  '
